# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Sun May 14 22:00:32 UTC 2023 with GitHub Actions"
# Updates the Price (column D) and Volume(1h) (column E) cells for rows 2-51.
# Price/Volume cells are stored as plain text in this sheet (not numbers), so
# number-looking Price values are written with a leading apostrophe to force
# Excel to keep them as literal text (preserves values like "85.50"/"1.040"
# instead of Excel normalising them to 85.5/1.04).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.775.98"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.857.11"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "'1.036"
$ws.Range("D5").Value = "'322.69"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'1.032"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").Value = "'0.4404"
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("D8").Value = "'0.3817"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").Value = "'0.07422"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").Value = "'0.8878"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "1.867.45"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "'5.525"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "'6.732"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "'0.07194"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "'85.50"
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("D17").Value = "'1.040"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").Value = "'0.000009084"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "'1.032"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "'15.55"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "27.777.05"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "'5.285"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "2.086.64"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").Value = "'2.058"
$ws.Range("E25").Value = "  +6.27%  "
$ws.Range("D26").Value = "'159.15"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").Value = "'18.75"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "'2.002"
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("D29").Value = "'5.354"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "'118.27"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("D31").Value = "'0.09109"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").Value = "'1.215"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").Value = "'0.7739"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("E34").Value = "  +4.71%  "
$ws.Range("D35").Value = "'4.601"
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("D36").Value = "'1.034"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").Value = "'1.157"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "'0.05295"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").Value = "'2.866"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").Value = "'0.5202"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").Value = "'6.911"
$ws.Range("E42").Value = "  +2.99%  "
$ws.Range("D43").Value = "'0.1674"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "'8.761"
$ws.Range("E44").Value = "  +2.51%  "
$ws.Range("D45").Value = "'110.71"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").Value = "'10.75"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").Value = "'0.06582"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("D49").Value = "'1.713"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "'0.4725"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("E51").Value = "  +0.53%  "
